$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.456.25'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '2.633.91'
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("E4").Value = '  -0.01%  '
$c = $ws.Range("D5")
$c.Value = "'581.94"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -2.40%  '
$c = $ws.Range("D6")
$c.Value = "'156.75"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.81%  '
$c = $ws.Range("D7")
$c.Value = "'0.645"
$c.Style = "Normal"
$ws.Range("E7").Value = '  +3.01%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -3.24%  '
$ws.Range("E10").Value = '  +0.31%  '
$c = $ws.Range("D11")
$c.Value = "'0.390"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.23%  '
$ws.Range("E12").Value = '  +0.00%  '
$c = $ws.Range("D13")
$c.Value = "'28.63"
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.43%  '
$c = $ws.Range("D14")
$c.Value = "'0.0000187"
$c.Style = "Normal"
$ws.Range("E14").Value = '  -4.31%  '
$ws.Range("D15").Value = '3.108.91'
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("D16").Value = '64.263.50'
$ws.Range("E16").Value = '  -1.72%  '
$ws.Range("D17").Value = '2.621.03'
$ws.Range("E17").Value = '  -1.39%  '
$c = $ws.Range("D18")
$c.Value = "'12.26"
$c.Style = "Normal"
$ws.Range("E18").Value = '  -2.88%  '
$c = $ws.Range("D19")
$c.Value = "'4.68"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.73%  '
$ws.Range("E20").Value = '  -0.33%  '
$c = $ws.Range("D21")
$c.Value = "'347.39"
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.37%  '
$c = $ws.Range("D22")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.01%  '
$c = $ws.Range("D23")
$c.Value = "'68.31"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c = $ws.Range("D24")
$c.Value = "'1.78"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +7.28%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Range("D25")
$c.Value = "'0.0000113"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.50%  '
$c = $ws.Range("D26")
$c.Value = "'9.45"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.76%  '
$c = $ws.Range("D27")
$c.Value = "'593.62"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +10.03%  '
$ws.Range("E28").Value = '  +0.90%  '
$c = $ws.Range("D29")
$c.Value = "'8.00"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.67%  '
$c = $ws.Range("D30")
$c.Value = "'0.161"
$c.Style = "Normal"
$ws.Range("E30").Value = '  -1.69%  '
$ws.Range("E31").Value = '  -0.01%  '
$c = $ws.Range("D32")
$c.Value = "'2.08"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.09%  '
$c = $ws.Range("D33")
$c.Value = "'6.70"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'1.73"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.56%  '
$c = $ws.Range("D35")
$c.Value = "'5.36"
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.40%  '
$ws.Range("E36").Value = '  -1.50%  '
$c = $ws.Range("D37")
$c.Value = "'20.02"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -1.79%  '
$c = $ws.Range("D38")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.07%  '
$c = $ws.Range("D39")
$c.Value = "'1.93"
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.13%  '
$c = $ws.Range("D40")
$c.Value = "'154.89"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.54%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  +6.21%  '
$c = $ws.Range("D43")
$c.Value = "'158.11"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.81%  '
$ws.Range("E44").Value = '  -1.49%  '
$c = $ws.Range("D45")
$c.Value = "'23.34"
$c.Style = "Normal"
$ws.Range("E45").Value = '  +3.82%  '
$c = $ws.Range("D46")
$c.Value = "'0.0603"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.43%  '
$c = $ws.Range("D47")
$c.Value = "'0.635"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("E48").Value = '  +2.78%  '
$c = $ws.Range("D50")
$c.Value = "'19.22"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -2.25%  '
$ws.Range("E51").Value = '  -5.91%  '
